$d = $word.ActiveDocument
$d.Paragraphs(1).Range.InsertAfter(" (Changed main)")
